$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Fibonacci recursion (depth 38) VFX (ms) value
$ws.Range("E17").Value = 37947

# Move the active selection to E18 as recorded in the saved view state
$ws.Range("E18").Select()
